# GoodInfo_v2 - 2021.11.19 raw data
#
# - Convert the plain "=O+P+Q" formulas in column R into three shared-formula
#   groups (R2:R33, R34:R65, R66:R97 - the third group's master ref
#   historically spans further than the live data, which only goes to row 68;
#   we recreate that by writing the formula across the full historical range
#   and then deleting the now-empty trailing rows so the shared-formula
#   "anchor" cell keeps the wider ref while the sheet's used range shrinks
#   back to row 68).
# - Hide the analysis input columns (D:M) that fed the summary columns.
# - Drop the now-unneeded helper column S (the "v"/marker column) entirely.
# - Leave the selection sitting on the (now empty) column S, matching the
#   state the workbook was left in after the cleanup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-enter the R column formulas as three shared-formula groups.
$ws.Range("R2:R33").Formula = "=O2+P2+Q2"
$ws.Range("R34:R65").Formula = "=O34+P34+Q34"
$ws.Range("R66:R97").Formula = "=O66+P66+Q66"

# Rows 69:97 were only ever touched by the R-column formula above (to give
# the third shared-formula group its historical R66:R97 ref) - remove them
# again so the sheet's real data still ends at row 68.
$ws.Range("A69:S97").EntireRow.Delete() | Out-Null

# Select column S before it disappears, mirroring the workbook's saved
# selection state (S1:S1048576).
$ws.Range("S:S").Select() | Out-Null

# Hide the raw data columns that back the summary/output columns.
$ws.Columns("D:E").Hidden = $true
$ws.Columns("F:L").Hidden = $true
$ws.Columns("M").Hidden = $true

# The helper/marker column S is no longer needed - delete it outright.
$ws.Columns("S").Delete() | Out-Null
